$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2, 4, 7, 8, 9, 10 with new odds values ---
$ws.Cells.Item(2, 6).Value = 2.76  # F2
$ws.Cells.Item(2, 8).Value = 2.52  # H2
$ws.Cells.Item(2, 12).Value = 1.31  # L2
$ws.Cells.Item(2, 15).Value = 1.22  # O2
$ws.Cells.Item(2, 16).Value = 2.32  # P2
$ws.Cells.Item(2, 17).Value = 1.67  # Q2
$ws.Cells.Item(2, 19).Value = 2.66  # S2
$ws.Cells.Item(2, 20).Value = 1.6  # T2
$ws.Cells.Item(2, 23).Value = 1.52  # W2
$ws.Cells.Item(2, 25).Value = 15  # Y2
$ws.Cells.Item(2, 28).Value = 16  # AB2
$ws.Cells.Item(2, 33).Value = 13.5  # AG2
$ws.Cells.Item(4, 6).Value = 1.68  # F4
$ws.Cells.Item(7, 13).Value = 1.05  # M7
$ws.Cells.Item(8, 14).Value = 1.1  # N8
$ws.Cells.Item(8, 16).Value = 2.52  # P8
$ws.Cells.Item(9, 15).Value = 1.51  # O9
$ws.Cells.Item(9, 17).Value = 2.44  # Q9
$ws.Cells.Item(9, 20).Value = 2.14  # T9
$ws.Cells.Item(9, 21).Value = 1.76  # U9
$ws.Cells.Item(10, 6).Value = 1.19  # F10
$ws.Cells.Item(10, 7).Value = 1.25  # G10
$ws.Cells.Item(10, 8).Value = 11  # H10
$ws.Cells.Item(10, 9).Value = 20  # I10
$ws.Cells.Item(10, 10).Value = 7  # J10
$ws.Cells.Item(10, 11).Value = 9.8  # K10
$ws.Cells.Item(10, 12).Value = 1.18  # L10
$ws.Cells.Item(10, 13).Value = 1.01  # M10
$ws.Cells.Item(10, 14).Value = 3.45  # N10
$ws.Cells.Item(10, 15).Value = 1.1  # O10
$ws.Cells.Item(10, 16).Value = 3.25  # P10
$ws.Cells.Item(10, 17).Value = 1.3  # Q10
$ws.Cells.Item(10, 18).Value = 2.04  # R10
$ws.Cells.Item(10, 19).Value = 1.71  # S10
$ws.Cells.Item(10, 20).Value = 1.8  # T10
$ws.Cells.Item(10, 21).Value = 1.98  # U10
$ws.Cells.Item(10, 22).Value = 1.06  # V10
$ws.Cells.Item(10, 23).Value = 4.8  # W10
$ws.Cells.Item(10, 24).Value = 55  # X10
$ws.Cells.Item(10, 25).Value = 75  # Y10
$ws.Cells.Item(10, 26).Value = 190  # Z10
$ws.Cells.Item(10, 27).Value = 610  # AA10
$ws.Cells.Item(10, 28).Value = 17.5  # AB10
$ws.Cells.Item(10, 29).Value = 23  # AC10
$ws.Cells.Item(10, 30).Value = 65  # AD10
$ws.Cells.Item(10, 31).Value = 230  # AE10
$ws.Cells.Item(10, 32).Value = 12.5  # AF10
$ws.Cells.Item(10, 33).Value = 14.5  # AG10
$ws.Cells.Item(10, 34).Value = 34  # AH10
$ws.Cells.Item(10, 35).Value = 160  # AI10
$ws.Cells.Item(10, 36).Value = 12  # AJ10
$ws.Cells.Item(10, 37).Value = 15  # AK10
$ws.Cells.Item(10, 38).Value = 34  # AL10
$ws.Cells.Item(10, 39).Value = 140  # AM10
$ws.Cells.Item(10, 40).Value = 3.05  # AN10
$ws.Cells.Item(10, 41).Value = 200  # AO10

# --- Insert a new row at position 11 for Dutch Eerste Divisie match ---
$ws.Rows.Item(11).Insert()

# --- Populate the newly inserted row 11 ---
# Column B holds a date-like string ("2025-12-04"). A direct .Value assignment
# would be auto-converted to a real date by Excel. Copy the equivalent text
# value from row 2 (values-only paste) so it stays plain text, like the rest
# of the column.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = "Dutch Eerste Divisie"  # A11
$ws.Cells.Item(11, 3).Value = "16:00:00"  # C11
$ws.Cells.Item(11, 4).Value = "Vitesse Arnhem"  # D11
$ws.Cells.Item(11, 5).Value = "De Graafschap"  # E11
$ws.Cells.Item(11, 6).Value = 2.68  # F11
$ws.Cells.Item(11, 7).Value = 3.3  # G11
$ws.Cells.Item(11, 8).Value = 2.26  # H11
$ws.Cells.Item(11, 9).Value = 2.56  # I11
$ws.Cells.Item(11, 10).Value = 3.95  # J11
$ws.Cells.Item(11, 11).Value = 4.7  # K11
$ws.Cells.Item(11, 12).Value = 1.01  # L11
$ws.Cells.Item(11, 13).Value = 1.02  # M11
$ws.Cells.Item(11, 14).Value = 2.8  # N11
$ws.Cells.Item(11, 15).Value = 1.14  # O11
$ws.Cells.Item(11, 16).Value = 2.8  # P11
$ws.Cells.Item(11, 17).Value = 1.46  # Q11
$ws.Cells.Item(11, 18).Value = 1.66  # R11
$ws.Cells.Item(11, 19).Value = 1.93  # S11
$ws.Cells.Item(11, 20).Value = 1.34  # T11
$ws.Cells.Item(11, 21).Value = 1.04  # U11
$ws.Cells.Item(11, 22).Value = 1.64  # V11
$ws.Cells.Item(11, 23).Value = 1.44  # W11
$ws.Cells.Item(11, 24).Value = 1000  # X11
$ws.Cells.Item(11, 25).Value = 1000  # Y11
$ws.Cells.Item(11, 26).Value = 1000  # Z11
$ws.Cells.Item(11, 27).Value = 1000  # AA11
$ws.Cells.Item(11, 28).Value = 1000  # AB11
$ws.Cells.Item(11, 29).Value = 1000  # AC11
$ws.Cells.Item(11, 30).Value = 1000  # AD11
$ws.Cells.Item(11, 31).Value = 1000  # AE11
$ws.Cells.Item(11, 32).Value = 1000  # AF11
$ws.Cells.Item(11, 33).Value = 1000  # AG11
$ws.Cells.Item(11, 34).Value = 1000  # AH11
$ws.Cells.Item(11, 35).Value = 1000  # AI11
$ws.Cells.Item(11, 36).Value = 1000  # AJ11
$ws.Cells.Item(11, 37).Value = 1000  # AK11
$ws.Cells.Item(11, 38).Value = 1000  # AL11
$ws.Cells.Item(11, 39).Value = 1000  # AM11
$ws.Cells.Item(11, 40).Value = 1000  # AN11
$ws.Cells.Item(11, 41).Value = 1000  # AO11

# --- Apply surgical updates to shifted rows 12, 14, 15 (row 13 unchanged) ---
$ws.Cells.Item(12, 9).Value = 7.6  # I12
$ws.Cells.Item(12, 12).Value = 1.26  # L12
$ws.Cells.Item(12, 21).Value = 2.24  # U12
$ws.Cells.Item(12, 22).Value = 1.15  # V12
$ws.Cells.Item(12, 23).Value = 3  # W12
$ws.Cells.Item(12, 24).Value = 27  # X12
$ws.Cells.Item(12, 27).Value = 220  # AA12
$ws.Cells.Item(12, 29).Value = 11.5  # AC12
$ws.Cells.Item(12, 30).Value = 26  # AD12
$ws.Cells.Item(12, 31).Value = 85  # AE12
$ws.Cells.Item(12, 32).Value = 10  # AF12
$ws.Cells.Item(12, 34).Value = 20  # AH12
$ws.Cells.Item(12, 35).Value = 75  # AI12
$ws.Cells.Item(12, 37).Value = 13  # AK12
$ws.Cells.Item(12, 38).Value = 26  # AL12
$ws.Cells.Item(12, 41).Value = 80  # AO12
$ws.Cells.Item(14, 6).Value = 2.06  # F14
$ws.Cells.Item(14, 7).Value = 2.18  # G14
$ws.Cells.Item(14, 8).Value = 3.8  # H14
$ws.Cells.Item(14, 10).Value = 3.45  # J14
$ws.Cells.Item(14, 14).Value = 3.2  # N14
$ws.Cells.Item(14, 15).Value = 1.38  # O14
$ws.Cells.Item(14, 16).Value = 1.76  # P14
$ws.Cells.Item(14, 17).Value = 2.06  # Q14
$ws.Cells.Item(14, 32).Value = 15.5  # AF14
$ws.Cells.Item(15, 17).Value = 1.83  # Q15
